$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header C1 text
$ws.Range("C1").Value = "desvio ($/200media)"

# Row 2 (AAPL) - update numeric values
$ws.Range("B2").Value = 0.9542463252951378
$ws.Range("C2").Value = 42.02133713945131
$ws.Range("D2").Value = -1.098849922273934

# Row 3: Ticker changes from IBM to MSFT, with new numeric values
$ws.Range("A3").Value = "MSFT"
$ws.Range("B3").Value = 1.175842304054456
$ws.Range("C3").Value = 70.2347054836018
$ws.Range("D3").Value = 0.8184572707958198

# Row 4 (AMZN) is removed entirely
$ws.Range("A4:D4").Delete()
